$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two trailing data rows (5 and 6) - shrinks the used range to A1:P4
$ws.Rows("5:6").Delete()

# New header row and data rows replacing the CSV-derived content
$data = @(
    @("Cert","Lender Loan #","City","State","Zip","Orig Loan Amt","Orig Appr Value","Orig Sales Price","Last Ins Amt","Renewal Period","Renewal Option","Refund Ind","Cov %","Loan Closing Date","Next Due Date","Days Past Due"),
    @("3470039372","329696-590","PHOENIX","AZ","85033","87400","93500","92000","87400","Zero Monthly","C","R","30","12/30/1998","01/01/2020","196"),
    @("3877345727","0579129166","WATSONTOWN","PA","17777","80000","80000","80000","80000","Zero Monthly","C","R","20","04/12/2007","05/01/2020","75"),
    @("6301710550","0579130324","SPOKANE","WA","99207","113400","114000","113400","113400","Zero Monthly","C","R","35","10/11/2007","05/01/2020","75")
)

# Force every cell in the range to Text format so numeric- and date-looking
# strings (e.g. "87400", "12/30/1998") are stored as text, matching the
# original inline-string ("t=inlineStr") CSV-import representation instead
# of being auto-converted to numbers/dates.
$ws.Range("A1:P4").NumberFormat = "@"

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}
